# Add Q11-Q18 (new geo-dimension requests) to "formalization of requirements"
# and their legend entries to "preliminary workload".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # formalization of requirements
$ws2 = $wb.Worksheets.Item(2)   # preliminary workload

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Formats first (so subsequent .Value writes land on cells that already
#    carry the correct style index - mirrors how the K column looks).
# ---------------------------------------------------------------------------

# Header row (row 1) L1:S1 -> same style as K1
$ws1.Range("K1").Copy()
$ws1.Range("L1:S1").PasteSpecial($xlPasteFormats)

# Body rows 3:17 and 19:24, columns L:S -> same style as K3 (row 18 is a
# section header with only column A populated, so it is skipped)
$ws1.Range("K3").Copy()
$ws1.Range("L3:S17").PasteSpecial($xlPasteFormats)
$ws1.Range("K3").Copy()
$ws1.Range("L19:S24").PasteSpecial($xlPasteFormats)

$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Values - written in the precise order the strings were first introduced
#    so the shared-string table indices line up with the target workbook.
# ---------------------------------------------------------------------------

# Q11 .. Q14
$ws1.Range("L1").Value = "Q11"
$ws1.Range("M1").Value = "Q12"
$ws1.Range("N1").Value = "Q13"
$ws1.Range("O1").Value = "Q14"

# River / building questions (legend sheet)
$ws2.Range("B12").Value = "What is the percentage of sick trees that are near a river ?"
$ws2.Range("B14").Value = "What is the percentage of sick trees that are near a building ?"
$ws2.Range("B15").Value = "How many trees are near a building ?"
$ws2.Range("B13").Value = "How many trees are near a river ?"

# Q15 .. Q18
$ws1.Range("P1").Value = "Q15"
$ws1.Range("Q1").Value = "Q16"
$ws1.Range("R1").Value = "Q17"
$ws1.Range("S1").Value = "Q18"

# Heliport / motorway questions (legend sheet)
$ws2.Range("B16").Value = "What is the percentage of sick trees that are near a heliport ?"
$ws2.Range("B17").Value = "How many trees are near a heliport ?"
$ws2.Range("B18").Value = "What is the percentage of sick trees that are near a motorway ?"
$ws2.Range("B19").Value = "How many trees are near a motorway ?"

# RequestId column on the legend sheet (reuses the Q11..Q18 strings above)
$ws2.Range("A12").Value = "Q11"
$ws2.Range("A13").Value = "Q12"
$ws2.Range("A14").Value = "Q13"
$ws2.Range("A15").Value = "Q14"
$ws2.Range("A16").Value = "Q15"
$ws2.Range("A17").Value = "Q16"
$ws2.Range("A18").Value = "Q17"
$ws2.Range("A19").Value = "Q18"

# ---------------------------------------------------------------------------
# 3) "x" marks on the requirement matrix
# ---------------------------------------------------------------------------

# Row 3 (species) is answered by every new request
$ws1.Range("L3:S3").Value = "x"

# Row 17 (isLast) - percentage questions (odd columns M,O,Q,S)
$ws1.Range("M17").Value = "x"
$ws1.Range("O17").Value = "x"
$ws1.Range("Q17").Value = "x"
$ws1.Range("S17").Value = "x"

# Row 19 (sickOnRoot) - percentage questions (odd columns M,O,Q,S)
$ws1.Range("M19").Value = "x"
$ws1.Range("O19").Value = "x"
$ws1.Range("Q19").Value = "x"
$ws1.Range("S19").Value = "x"

# Row 20 (quantity) - count questions (even columns L,N,P,R)
$ws1.Range("L20").Value = "x"
$ws1.Range("N20").Value = "x"
$ws1.Range("P20").Value = "x"
$ws1.Range("R20").Value = "x"

# ---------------------------------------------------------------------------
# 4) Column widths + selections to match the authored state
# ---------------------------------------------------------------------------

$ws1.Columns("L:S").ColumnWidth = $ws1.Columns("K:K").ColumnWidth

$ws1.Range("L19").Select()
$ws2.Range("B20").Select()
